$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header-ish cells B64:B67 (plain shared-string values A/B/C/D) ---
$ws.Range("B64").Value = "A"
$ws.Range("B65").Value = "B"
$ws.Range("B66").Value = "C"
$ws.Range("B67").Value = "D"

# --- B69:B80 : broadcast first column of _nData across its width, stacked ---
$f1 = '=LET(_c,COLUMNS(_nData)-1,DROP(REDUCE("",TAKE(_nData,,1),LAMBDA(a,v,VSTACK(a,EXPAND(v,_c,1,v)))),1))'
$ws.Range("B69:B80").FormulaArray = $f1

# --- C69:C80 : remaining columns of _nData flattened to one column ---
$f2 = '=TOCOL(DROP(_nData,,1))'
$ws.Range("C69:C80").FormulaArray = $f2

# --- F69:G80 : same broadcast logic wrapped in a LAMBDA, applied to _nData ---
$f3 = "=LAMBDA(_nData,LET(f, LAMBDA(_d,LET(_c,COLUMNS(_d)-1,DROP(REDUCE(`"`",TAKE(_d,,1),LAMBDA(a,v,VSTACK(a,EXPAND(v,_c,1,v)))),1))),
     HSTACK(f(_nData),TOCOL(DROP(_nData,,1)))
))(_nData)"
$ws.Range("F69:G80").FormulaArray = $f3

# --- I69:J83 : same LAMBDA, applied to _nAlt instead ---
$f4 = "=LAMBDA(_nData,LET(f, LAMBDA(_d,LET(_c,COLUMNS(_d)-1,DROP(REDUCE(`"`",TAKE(_d,,1),LAMBDA(a,v,VSTACK(a,EXPAND(v,_c,1,v)))),1))),
     HSTACK(f(_nData),TOCOL(DROP(_nData,,1)))
))(_nAlt)"
$ws.Range("I69:J83").FormulaArray = $f4

# --- view state: scroll position + active selection ---
$excel.ActiveWindow.ScrollRow = 48
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("K65").Select()
